# Rename the "congenital" variable name to "misc_long_term" across every
# variables_* worksheet in the workbook. Each worksheet stores its list of
# variable names as inline strings in column A (rows 2-6); exactly one of
# those rows currently reads "congenital" on the sheets that have it.
$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

$oldName = "congenital"
$newName = "misc_long_term"
$replacedCount = 0

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $usedRange = $ws.UsedRange

    foreach ($cell in $usedRange.Cells) {
        if ($cell.Text -eq $oldName) {
            $cell.Value = $newName
            $replacedCount = $replacedCount + 1
        }
    }
}

Write-Output "Replaced '$oldName' with '$newName' in $replacedCount cell(s)"
